# Apply the "crypto price snapshot" refresh described in the commit.
# Price cells that look like plain decimals (e.g. "0.999", "1.00", "17.93")
# are written with a leading apostrophe so Excel keeps them as literal text
# instead of re-parsing/rounding them as numbers (matching how the sheet
# originally stored every Price/Volume column as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.464.98'
$ws.Range("D3").Value = '3.744.65'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''594.90'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = '''167.55'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("D7").Value = '3.745.98'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("D11").Value = '''6.48'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E13").Value = '  -5.90%  '
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '4.372.90'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '3.743.52'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = '68.431.90'
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").Value = '''17.93'
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").Value = '''467.58'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '''0.701'
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("D24").Value = '''84.17'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("D26").Value = '''2.20'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '''10.14'
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D30").Value = '3.890.37'
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("D32").Value = '''7.34'
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("D33").Value = '''29.95'
$ws.Range("E33").Value = '  -1.58%  '
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("D35").Value = '''9.28'
$ws.Range("E35").Value = '  +1.90%  '
$ws.Range("D36").Value = '''0.999'
$ws.Range("D37").Value = '3.699.76'
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").Value = '''3.44'
$ws.Range("E39").Value = '  -9.71%  '
$ws.Range("D40").Value = '''0.139'
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '''1.94'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value = '''8.62'
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("D48").Value = '''42.84'
$ws.Range("E48").Value = '  +9.91%  '
$ws.Range("D49").Value = '''45.88'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").Value = '''146.94'
$ws.Range("E50").Value = '  +5.78%  '
$ws.Range("D51").Value = '''392.88'
$ws.Range("E51").Value = '  -1.20%  '
